# Update the "ballots" worksheet with the 12/17 ballots (rows 50-53)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# Row 50 - Roberto Colon (source: Unanimo Sports)
$ws.Range("A50").Value = "Roberto Colon"
$ws.Range("C50").Value = "x"
$ws.Range("D50").Value = "x"
$ws.Range("E50").Value = "x"
$ws.Range("I50").Value = "x"
$ws.Range("J50").Value = "x"
$ws.Range("K50").Value = "x"
$ws.Range("O50").Value = "x"
$ws.Range("Q50").Value = "x"
$ws.Range("S50").Value = "x"
$ws.Range("T50").Value = "x"
$ws.Range("AK50").Value = 10
$ws.Range("AL50").Value = "Unanimo Sports"
$ws.Range("AM50").Value = 43451

# Row 51 - Greg Cote (source: Twitter)
$ws.Range("A51").Value = "Greg Cote"
$ws.Range("C51").Value = "x"
$ws.Range("D51").Value = "x"
$ws.Range("E51").Value = "x"
$ws.Range("I51").Value = "x"
$ws.Range("J51").Value = "x"
$ws.Range("K51").Value = "x"
$ws.Range("N51").Value = "x"
$ws.Range("O51").Value = "x"
$ws.Range("R51").Value = "x"
$ws.Range("AK51").Value = 9
$ws.Range("AL51").Value = "Twitter"
$ws.Range("AM51").Value = 43451

# Row 52 - Felix DeJesus (source: Twitter)
$ws.Range("A52").Value = "Felix DeJesus"
$ws.Range("C52").Value = "x"
$ws.Range("D52").Value = "x"
$ws.Range("E52").Value = "x"
$ws.Range("I52").Value = "x"
$ws.Range("K52").Value = "x"
$ws.Range("N52").Value = "x"
$ws.Range("O52").Value = "x"
$ws.Range("Q52").Value = "x"
$ws.Range("S52").Value = "x"
$ws.Range("T52").Value = "x"
$ws.Range("AK52").Value = 10
$ws.Range("AL52").Value = "Twitter"
$ws.Range("AM52").Value = 43451

# Row 53 - Mike Vaccaro (source: Twitter)
$ws.Range("A53").Value = "Mike Vaccaro"
$ws.Range("C53").Value = "x"
$ws.Range("D53").Value = "x"
$ws.Range("E53").Value = "x"
$ws.Range("I53").Value = "x"
$ws.Range("J53").Value = "x"
$ws.Range("K53").Value = "x"
$ws.Range("O53").Value = "x"
$ws.Range("Q53").Value = "x"
$ws.Range("U53").Value = "x"
$ws.Range("V53").Value = "x"
$ws.Range("AK53").Value = 10
$ws.Range("AL53").Value = "Twitter"
$ws.Range("AM53").Value = 43451

# Copy the date number format from the existing AM49 cell onto the new
# date cells so they reuse the same cell style (instead of creating a
# new numFmt/style entry).
$ws.Range("AM49").Copy()
$ws.Range("AM50:AM53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move/select the new last cell, mirroring where the user ended up after
# entering the new ballots.
$ws.Range("A53").Select()
